# Rename the "HO" sheet to "YOY OF HO", make it the active/selected tab,
# and move its selection to cell C19 (previously the active sheet was the
# first one, with this sheet's selection sitting on G1:G2).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("HO")
$ws.Name = "YOY OF HO"

# Activating the sheet makes it the workbook's active tab (activeTab) and
# marks its own sheetView as tabSelected, while the previously-active first
# sheet naturally drops tabSelected.
$ws.Activate()

# Move the selection on this sheet from G1:G2 to C19.
$ws.Range("C19").Select()
